# Fix for "pregunta 2" (question 2): the "materia" column (H) was being
# included in this "solo desconocidos" (unknown-only) teacher segregation,
# which broke the report. Remove the entire "materia" column so that
# "edad" and "id_centro" shift left into H/I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(8).Delete()
